# Fill in the evaluation grades on the "Projet" sheet (first sheet of the workbook).
# The workbook has 4 evaluation blocks (columns B/C/D..E, H/I/J..K, N/O/P..Q, T/U/V..W)
# each holding scored rows. This commit fills in the previously-empty "score" column
# for each criterion, which in turn drives all the dependent formulas/totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Projet")

$ws.Range("B5").Value = 1
$ws.Range("H5").Value = 2
$ws.Range("N5").Value = 2
$ws.Range("T5").Value = 1

$ws.Range("H6").Value = 1
$ws.Range("N6").Value = 1
$ws.Range("T6").Value = 1

$ws.Range("B7").Value = 1
$ws.Range("H7").Value = 2
$ws.Range("N7").Value = 1
$ws.Range("T7").Value = 1

$ws.Range("B8").Value = 1
$ws.Range("H8").Value = 2
$ws.Range("N8").Value = 1
$ws.Range("T8").Value = 1

$ws.Range("H9").Value = 2
$ws.Range("N9").Value = 2
$ws.Range("T9").Value = 1

$ws.Range("B10").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("N10").Value = 2
$ws.Range("T10").Value = 1

$ws.Range("H11").Value = 1
$ws.Range("N11").Value = 2
$ws.Range("T11").Value = 1

$ws.Range("H12").Value = 1
$ws.Range("N12").Value = 1
$ws.Range("T12").Value = 1

$ws.Range("H13").Value = 1
$ws.Range("N13").Value = 2
$ws.Range("T13").Value = 1

$ws.Range("B14").Value = 1
$ws.Range("H14").Value = 2
$ws.Range("N14").Value = 3
$ws.Range("T14").Value = 1

$ws.Range("B15").Value = 1
$ws.Range("H15").Value = 2

$ws.Range("B16").Value = 1
$ws.Range("H16").Value = 2

$ws.Range("B17").Value = 1
$ws.Range("H17").Value = 1

$ws.Range("H18").Value = 2

$ws.Range("H19").Value = 2

$ws.Range("H20").Value = 2

# Move the active selection to B15, matching the author's last cursor position.
$ws.Range("B15").Select()
